$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new column L, mirroring column K's formatting ---
# Copy formatting (styles) from K4:K12 down to L4:L12 so the new cells
# pick up the same number formats / borders / fonts as their K counterparts.
$ws.Range("K4:K12").Copy()
$ws.Range("L4:L12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 4 header gets the new year 2020
$ws.Range("L4").Value = 2020

# Rows 5-12 duplicate the values already present in column K
$ws.Range("L5").Value  = 5.6
$ws.Range("L6").Value  = 0.8
$ws.Range("L7").Value  = 1.9
$ws.Range("L8").Value  = 0.7
$ws.Range("L9").Value  = 0.7
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.3
$ws.Range("L12").Value = 0.2

# --- Update the current selection shown in the workbook ---
$null = $ws.Range("N5").Select()
